$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data through row 103. We append 4 new workout log
# rows (104-107). Seed the new rows by copying the last existing data row so
# that formatting/styles (e.g. the date style on column B) are carried over
# exactly, then overwrite the individual cell values with the real data.

$ws.Range("A103:M103").Copy($ws.Range("A104:M104"))
$ws.Range("A103:M103").Copy($ws.Range("A105:M105"))
$ws.Range("A103:M103").Copy($ws.Range("A106:M106"))
$ws.Range("A103:M103").Copy($ws.Range("A107:M107"))

# Row 104: Jeremiah / Run / Agile Antelope
$ws.Range("A104").Value = "Jeremiah"
$ws.Range("B104").Value = 45467
$ws.Range("C104").Value = "Run"
$ws.Range("D104").Value = 10
$ws.Range("E104").Value = 1.01
$ws.Range("F104").Value = 89
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 2
$ws.Range("I104").Value = 4
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = "Agile Antelope"
$ws.Range("M104").Value = 3

# Row 105: Jeremiah / Workout / Agile Antelope
$ws.Range("A105").Value = "Jeremiah"
$ws.Range("B105").Value = 45467
$ws.Range("C105").Value = "Workout"
$ws.Range("D105").Value = 41
$ws.Range("E105").Value = 0
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 11
$ws.Range("H105").Value = 22
$ws.Range("I105").Value = 8
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = "Agile Antelope"
$ws.Range("M105").Value = 3

# Row 106: Matt / Ride / Agile Antelope
$ws.Range("A106").Value = "Matt"
$ws.Range("B106").Value = 45467
$ws.Range("C106").Value = "Ride"
$ws.Range("D106").Value = 32
$ws.Range("E106").Value = 10.68
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 28
$ws.Range("I106").Value = 3
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = "Agile Antelope"
$ws.Range("M106").Value = 3

# Row 107: Steven / Walk / Wily Hyena (new workout-level string)
$ws.Range("A107").Value = "Steven"
$ws.Range("B107").Value = 45467
$ws.Range("C107").Value = "Walk"
$ws.Range("D107").Value = 28
$ws.Range("E107").Value = 1.4
$ws.Range("F107").Value = 62
$ws.Range("G107").Value = 28
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = "Wily Hyena"
$ws.Range("M107").Value = 3

# Update the view so the active cell matches where the author ended up after
# entering the new rows.
$ws.Range("E108").Select()
